# Updating filtered feeds from workflow
# Appends two new rows to the "Filtered Feeds" sheet for the Promega
# OncoMate MSI Dx Analysis System CDx story (one row per source link).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$genomewebLink = "https://www.genomeweb.com/cancer/promega-receives-fda-approval-oncomate-msi-dx-analysis-system-endometrial-cancer-cdx"
$dx360Link     = "https://www.360dx.com/cancer/promega-receives-fda-approval-oncomate-msi-dx-analysis-system-endometrial-cancer-cdx"
$keywords      = "CDx"
$title         = "Promega Receives FDA Approval for OncoMate MSI Dx Analysis System as Endometrial Cancer CDx"

# Row 61 - GenomeWeb source
$ws.Range("A61").Value = $genomewebLink
$ws.Range("B61").Value = $keywords
$ws.Range("C61").Value = $title

# Row 62 - 360Dx source
$ws.Range("A62").Value = $dx360Link
$ws.Range("B62").Value = $keywords
$ws.Range("C62").Value = $title

# Turn column A entries into real hyperlinks (adds the relationship + rId)
$ws.Hyperlinks.Add($ws.Range("A61"), $genomewebLink)
$ws.Hyperlinks.Add($ws.Range("A62"), $dx360Link)

# Match the hyperlink cell styling already used by the rest of column A
$ws.Range("A61").Style = $ws.Range("A2").Style
$ws.Range("A62").Style = $ws.Range("A2").Style
